$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: apply the identical timetable restructuring to a Section sheet
# (Section_A / Section_B). Both sheets receive exactly the same edits.
# ---------------------------------------------------------------------------
function Update-SectionSheet($ws) {
    # Row 2 — 09:00-10:30 slot
    $ws.Range("A2").Value = "09:00-10:30"
    $ws.Range("B2").Value = "DS456 (Elective)"
    $ws.Range("C2").Value = "Free"
    $ws.Range("D2").Value = "EC456 (Elective)"
    $ws.Range("E2").Value = "DS456 (Elective)"
    $ws.Range("F2").Value = "Free"

    # Row 3 — 10:30-12:00 slot (time label unchanged)
    $ws.Range("B3").Value = "Free"
    $ws.Range("C3").Value = "DS401 (Elective)"
    $ws.Range("D3").Value = "Free"
    $ws.Range("E3").Value = "Free"
    $ws.Range("F3").Value = "Free"

    # Row 4 — lunch break, time label changes
    $ws.Range("A4").Value = "12:00-13:00"

    # Row 5 — time label changes
    $ws.Range("A5").Value = "13:00-14:30"
    $ws.Range("B5").Value = "Free"
    $ws.Range("C5").Value = "Free"
    $ws.Range("D5").Value = "EC456 (Elective)"
    $ws.Range("E5").Value = "Free"
    $ws.Range("F5").Value = "Free"

    # Row 6 — time label changes
    $ws.Range("A6").Value = "14:30-15:30"
    $ws.Range("B6").Value = "Free"
    $ws.Range("C6").Value = "Free"
    $ws.Range("D6").Value = "DS456 (Tutorial)"
    $ws.Range("E6").Value = "Free"
    $ws.Range("F6").Value = "Free"

    # Row 7 — time label changes
    $ws.Range("A7").Value = "15:30-17:00"
    $ws.Range("B7").Value = "Free"
    $ws.Range("C7").Value = "Free"
    $ws.Range("D7").Value = "Free"
    $ws.Range("E7").Value = "Free"
    $ws.Range("F7").Value = "DS401 (Elective)"

    # Row 8 — time label changes
    $ws.Range("A8").Value = "17:00-18:00"
    $ws.Range("B8").Value = "Free"
    $ws.Range("C8").Value = "Free"
    $ws.Range("D8").Value = "Free"
    $ws.Range("E8").Value = "DS401 (Tutorial)"
    $ws.Range("F8").Value = "EC456 (Tutorial)"

    # Rows 9-12 are no longer part of the timetable; remove them so the
    # sheet's used range (dimension) shrinks back down to A1:F8.
    $ws.Rows("9:12").Delete()
}

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")
Update-SectionSheet $wsA
Update-SectionSheet $wsB

# ---------------------------------------------------------------------------
# Course_Summary — swap in the new instructor names
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Course_Summary")
$wsSummary.Range("H2").Value = "Dr. Rajendra Hegadi"
$wsSummary.Range("H3").Value = "Dr. Divyajyothi"
$wsSummary.Range("H4").Value = "Dr. Girish G N"

# ---------------------------------------------------------------------------
# Elective_Coordination — restructure into the LTPSC-compliant layout:
#   Elective Course | Session Type | Day | Time Slot | Duration | Sections
# The old layout (Elective Course | Day | Time Slot | Slot Name | Sections)
# does not map cleanly column-for-column (the "Slot Name" column is dropped
# entirely and replaced by new per-row "Duration" data), so the simplest,
# most reliable approach is to wipe the existing values and rewrite the
# whole table, rather than trying to insert/shift columns in place.
# ---------------------------------------------------------------------------
$wsElec = $wb.Worksheets.Item("Elective_Coordination")

# Drop all existing data (ClearContents keeps the header row's bold/border
# style on A1:E1 so we don't have to rebuild it from scratch).
$wsElec.Cells.ClearContents()

# Column F ("Sections") does not exist yet — stamp it with the same header
# style as the neighbouring styled header cell (E1) before filling it in.
$wsElec.Range("E1").Copy()
$wsElec.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wsElec.Range("A1").Value = "Elective Course"
$wsElec.Range("B1").Value = "Session Type"
$wsElec.Range("C1").Value = "Day"
$wsElec.Range("D1").Value = "Time Slot"
$wsElec.Range("E1").Value = "Duration"
$wsElec.Range("F1").Value = "Sections"

$rows = @(
    @("DS456", "Lecture 1", "Mon", "09:00-10:30", "1.5 hours", "A & B (Common Slot)"),
    @("DS456", "Lecture 2", "Thu", "09:00-10:30", "1.5 hours", "A & B (Common Slot)"),
    @("DS456", "Tutorial",  "Wed", "14:30-15:30", "1 hour",    "A & B (Common Slot)"),
    @("EC456", "Lecture 1", "Wed", "13:00-14:30", "1.5 hours", "A & B (Common Slot)"),
    @("EC456", "Lecture 2", "Wed", "09:00-10:30", "1.5 hours", "A & B (Common Slot)"),
    @("EC456", "Tutorial",  "Fri", "17:00-18:00", "1 hour",    "A & B (Common Slot)"),
    @("DS401", "Lecture 1", "Fri", "15:30-17:00", "1.5 hours", "A & B (Common Slot)"),
    @("DS401", "Lecture 2", "Tue", "10:30-12:00", "1.5 hours", "A & B (Common Slot)"),
    @("DS401", "Tutorial",  "Thu", "17:00-18:00", "1 hour",    "A & B (Common Slot)")
)

$r = 2
foreach ($row in $rows) {
    $wsElec.Range("A$r").Value = $row[0]
    $wsElec.Range("B$r").Value = $row[1]
    $wsElec.Range("C$r").Value = $row[2]
    $wsElec.Range("D$r").Value = $row[3]
    $wsElec.Range("E$r").Value = $row[4]
    $wsElec.Range("F$r").Value = $row[5]
    $r++
}
